$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6498
$ws1.Range("F3").Value = 34
$ws1.Range("F6").Value = 119

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6498
$ws4.Range("F3").Value = 34
$ws4.Range("F6").Value = 119
